$d = $word.ActiveDocument

# --- 1) First paragraph: append trailing spaces to the existing run, then
#        append three new red-colored runs forming
#        "(This is a change – Version for branch alternate)"
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# Trim the paragraph mark off the end of the paragraph range so inserts land
# right after the existing text, before the pilcrow.
$r1.End = $r1.End - 1

# Append two trailing spaces to the original sentence (same run/sentence).
$r1.InsertAfter("  ")

# Run 2: "(This is a change – Ve"
$run2 = $d.Range($r1.End, $r1.End)
$run2.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run2.Font.Color = 192  # 0xC00000 as a BGR-ordered COM color long

# Run 3: "rsion for branch alternate"
$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter("rsion for branch alternate")
$run3.Font.Color = 192

# Run 4: ")"
$run4 = $d.Range($run3.End, $run3.End)
$run4.InsertAfter(")")
$run4.Font.Color = 192

# --- 2) Add a brand-new, plain empty paragraph (shaded F9F9F9) right after
#        the final "Free at last..." paragraph, before the section break.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$lastRange = $lastPara.Range
$lastRange.Collapse(0)  # wdCollapseEnd
$lastRange.InsertParagraphAfter() | Out-Null

$addedPara = $d.Paragraphs($d.Paragraphs.Count)
$addedXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$addedPara.Range.InsertXML($addedXml)
